$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 27.37907883418829
$ws.Range("C2").Value = 16.52175977624685
$ws.Range("D2").Value = 6.053805848669706
$ws.Range("E2").Value = 7.326001658998125
$ws.Range("G2").Value = 3.742306077781016
$ws.Range("I2").Value = 36.47353691899758
$ws.Range("L2").Value = 11.24777033414911
$ws.Range("N2").Value = 21.25543154662104

$ws.Range("B3").Value = 26.88870003491815
$ws.Range("C3").Value = 15.89276050444071
$ws.Range("D3").Value = 5.949053628159719
$ws.Range("E3").Value = 7.306207804260249
$ws.Range("G3").Value = 3.74808147660402
$ws.Range("I3").Value = 36.30476402954628
$ws.Range("L3").Value = 11.22912619536125
$ws.Range("N3").Value = 21.28723491231438

$ws.Range("B4").Value = 26.59295786546052
$ws.Range("C4").Value = 15.5001022417505
$ws.Range("D4").Value = 5.885970466723176
$ws.Range("E4").Value = 7.293845820229543
$ws.Range("G4").Value = 3.751802919987089
$ws.Range("I4").Value = 36.20980109898311
$ws.Range("L4").Value = 11.22033054298903
$ws.Range("N4").Value = 21.30868901273648

$ws.Range("B5").Value = 26.473954012554
$ws.Range("C5").Value = 15.33874859849085
$ws.Range("D5").Value = 5.8606081930798
$ws.Range("E5").Value = 7.288755202917619
$ws.Range("G5").Value = 3.753363742165789
$ws.Range("I5").Value = 36.17329430043787
$ws.Range("L5").Value = 11.21741351900184
$ws.Range("N5").Value = 21.31791422662682

$ws.Range("B6").Value = 26.45428977765084
$ws.Range("C6").Value = 15.31188355455339
$ws.Range("D6").Value = 5.856418604009487
$ws.Range("E6").Value = 7.287906691066093
$ws.Range("G6").Value = 3.753625597687563
$ws.Range("I6").Value = 36.16736503597656
$ws.Range("L6").Value = 11.21696945317779
$ws.Range("N6").Value = 21.31947514997549

$ws.Range("B7").Value = 26.59134658997167
$ws.Range("C7").Value = 15.49793121631057
$ws.Range("D7").Value = 5.885626982312188
$ws.Range("E7").Value = 7.293777381790896
$ws.Range("G7").Value = 3.751823790107389
$ws.Range("I7").Value = 36.20929986831292
$ws.Range("L7").Value = 11.22028850102365
$ws.Range("N7").Value = 21.30881147619865

$ws.Range("B8").Value = 27.20898640010407
$ws.Range("C8").Value = 16.30637947561581
$ws.Range("D8").Value = 6.017451519207547
$ws.Range("E8").Value = 7.319219651519279
$ws.Range("G8").Value = 3.744261177594781
$ws.Range("I8").Value = 36.41355426206934
$ws.Range("L8").Value = 11.24079183984288
$ws.Range("N8").Value = 21.26599625525778

$ws.Range("B9").Value = 28.45523255496243
$ws.Range("C9").Value = 17.82958253783069
$ws.Range("D9").Value = 6.2842446664792
$ws.Range("E9").Value = 7.367485578775018
$ws.Range("G9").Value = 3.730812127742364
$ws.Range("I9").Value = 36.88215322692703
$ws.Range("L9").Value = 11.30200040673324
$ws.Range("N9").Value = 21.19741385447108

$ws.Range("B10").Value = 29.38236676691933
$ws.Range("C10").Value = 18.89788789808266
$ws.Range("D10").Value = 6.483283536285496
$ws.Range("E10").Value = 7.401995962729109
$ws.Range("G10").Value = 3.721759213915345
$ws.Range("I10").Value = 37.26681127270234
$ws.Range("L10").Value = 11.35968865144121
$ws.Range("N10").Value = 21.15652776607866

$ws.Range("B11").Value = 29.80461894059443
$ws.Range("C11").Value = 19.37062654542875
$ws.Range("D11").Value = 6.574068763846064
$ws.Range("E11").Value = 7.417498158711773
$ws.Range("G11").Value = 3.71781762032989
$ws.Range("I11").Value = 37.45027509866243
$ws.Range("L11").Value = 11.38866504487844
$ws.Range("N11").Value = 21.14001860899775

$ws.Range("B12").Value = 29.96442402376906
$ws.Range("C12").Value = 19.54758113357756
$ws.Range("D12").Value = 6.608447716854101
$ws.Range("E12").Value = 7.423340816699546
$ws.Range("G12").Value = 3.716350203872723
$ws.Range("I12").Value = 37.52093843908752
$ws.Range("L12").Value = 11.40002735514644
$ws.Range("N12").Value = 21.13406992897497

$ws.Range("B13").Value = 29.93001379381067
$ws.Range("C13").Value = 19.50956476190869
$ws.Range("D13").Value = 6.601044099606453
$ws.Range("E13").Value = 7.422083717476349
$ws.Range("G13").Value = 3.716665121885523
$ws.Range("I13").Value = 37.50566741996237
$ws.Range("L13").Value = 11.39756301336961
$ws.Range("N13").Value = 21.13533756956118

$ws.Range("B14").Value = 29.81776897115556
$ws.Range("C14").Value = 19.38522684574462
$ws.Range("D14").Value = 6.57689732091187
$ws.Range("E14").Value = 7.417979389099987
$ws.Range("G14").Value = 3.717696391716432
$ws.Range("I14").Value = 37.45606495350355
$ws.Range("L14").Value = 11.38959204328335
$ws.Range("N14").Value = 21.13952312137141

$ws.Range("B15").Value = 29.74899896812827
$ws.Range("C15").Value = 19.30879352161244
$ws.Range("D15").Value = 6.56210580233996
$ws.Range("E15").Value = 7.415461781849094
$ws.Range("G15").Value = 3.718331346874767
$ws.Range("I15").Value = 37.42583599786229
$ws.Range("L15").Value = 11.38476021509472
$ws.Range("N15").Value = 21.14212642200981

$ws.Range("B16").Value = 29.35476739088118
$ws.Range("C16").Value = 18.86671296944215
$ws.Range("D16").Value = 6.477352417872405
$ws.Range("E16").Value = 7.400978950255019
$ws.Range("G16").Value = 3.722020341751744
$ws.Range("I16").Value = 37.25498963034338
$ws.Range("L16").Value = 11.35784969587935
$ws.Range("N16").Value = 21.15764894323971

$ws.Range("B17").Value = 29.11292708585746
$ws.Range("C17").Value = 18.59200126874419
$ws.Range("D17").Value = 6.425396106861422
$ws.Range("E17").Value = 7.392044238579293
$ws.Range("G17").Value = 3.724328502562647
$ws.Range("I17").Value = 37.15233252716974
$ws.Range("L17").Value = 11.34203879259628
$ws.Range("N17").Value = 21.1677085511223

$ws.Range("B18").Value = 28.97388285704452
$ws.Range("C18").Value = 18.43275628633413
$ws.Range("D18").Value = 6.395536909454459
$ws.Range("E18").Value = 7.386886673103649
$ws.Range("G18").Value = 3.725672733264676
$ws.Range("I18").Value = 37.09408709745734
$ws.Range("L18").Value = 11.33320233493481
$ws.Range("N18").Value = 21.17369115756948

$ws.Range("B19").Value = 28.92681957120779
$ws.Range("C19").Value = 18.37863117670703
$ws.Range("D19").Value = 6.385432401487392
$ws.Range("E19").Value = 7.385137204782553
$ws.Range("G19").Value = 3.726130730858658
$ws.Range("I19").Value = 37.07450450838001
$ws.Range("L19").Value = 11.33025479651837
$ws.Range("N19").Value = 21.17575045585382

$ws.Range("B20").Value = 29.13866681919681
$ws.Range("C20").Value = 18.62137414230729
$ws.Range("D20").Value = 6.430924640133745
$ws.Range("E20").Value = 7.392997269957625
$ws.Range("G20").Value = 3.724081074551877
$ws.Range("I20").Value = 37.16317794088265
$ws.Range("L20").Value = 11.34369525810102
$ws.Range("N20").Value = 21.1666173253655

$ws.Range("B21").Value = 29.85074178159982
$ws.Range("C21").Value = 19.42180500113345
$ws.Range("D21").Value = 6.583990065880825
$ws.Range("E21").Value = 7.419185676407398
$ws.Range("G21").Value = 3.717392801150418
$ws.Range("I21").Value = 37.47060236617362
$ws.Range("L21").Value = 11.39192276803355
$ws.Range("N21").Value = 21.1382854810834

$ws.Range("B22").Value = 30.31551298390173
$ws.Range("C22").Value = 19.93285931368951
$ws.Range("D22").Value = 6.684015525653033
$ws.Range("E22").Value = 7.436140217309253
$ws.Range("G22").Value = 3.713168291096339
$ws.Range("I22").Value = 37.67844021815875
$ws.Range("L22").Value = 11.42571100031908
$ws.Range("N22").Value = 21.12153612793239

$ws.Range("B23").Value = 30.06756401878071
$ws.Range("C23").Value = 19.66125151268365
$ws.Range("D23").Value = 6.630642064650854
$ws.Range("E23").Value = 7.427105798507072
$ws.Range("G23").Value = 3.715409644366605
$ws.Range("I23").Value = 37.56689072627433
$ws.Range("L23").Value = 11.40747129064644
$ws.Range("N23").Value = 21.13031306969693

$ws.Range("B24").Value = 29.12702989511072
$ws.Range("C24").Value = 18.60809873300056
$ws.Range("D24").Value = 6.428425152818375
$ws.Range("E24").Value = 7.392566469555765
$ws.Range("G24").Value = 3.724192883045671
$ws.Range("I24").Value = 37.15827231779794
$ws.Range("L24").Value = 11.34294558037205
$ws.Range("N24").Value = 21.16711004823764

$ws.Range("B25").Value = 28.11536409495559
$ws.Range("C25").Value = 17.42559936054947
$ws.Range("D25").Value = 6.211393463836417
$ws.Range("E25").Value = 7.354598258270794
$ws.Range("G25").Value = 3.734304020701252
$ws.Range("I25").Value = 36.74820202118333
$ws.Range("L25").Value = 11.28320008861972
$ws.Range("N25").Value = 21.21430792684682
